# Teaching overview: add the 2021 "Statistics in Connected Healthcare" course
# next to the existing 2020 entry in the "year" column (row 12, "Statistics
# in Connected Healthcare" topic).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("B12")
$cell.Value = "[2020](https://open.hpi.de/courses/StatisticsCHealthSoSe2020) + " + [char]10 + "[2021](https://open.hpi.de/courses/hpi-dh-StatisticsCHealth2021)"
